# Applies the diff: adds yellow highlighting to the "Exercicio 1" block on
# slide 11, renumbers/promotes the exercises on slide 12 (new "Exercicio 6"
# highlighted block inserted, old 1-5 shifted to 7-10), and renumbers the
# exercises on slide 13 (6-8 -> 11-13).

function ConvertTo-ComRgb {
    # PowerPoint/VBA RGB longs are packed as 0x00BBGGRR, i.e. the reverse of
    # the usual "RRGGBB" hex notation used in OOXML's <a:srgbClr val="..."/>.
    param([string]$RRGGBB)
    $r = [Convert]::ToInt32($RRGGBB.Substring(0,2), 16)
    $g = [Convert]::ToInt32($RRGGBB.Substring(2,2), 16)
    $b = [Convert]::ToInt32($RRGGBB.Substring(4,2), 16)
    return $r -bor ($g * 256) -bor ($b * 65536)
}

$yellow = ConvertTo-ComRgb "FFFF00"

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 11: highlight the title + description of "Exercicio 1" (the first
# two paragraphs of the content placeholder) in yellow. Nothing else on
# this slide changes.
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$body11 = $s11.Shapes.Item(2).TextFrame2.TextRange
$body11.Paragraphs(1).Font.Highlight.RGB = $yellow
$body11.Paragraphs(2).Font.Highlight.RGB = $yellow

# ---------------------------------------------------------------------
# Slide 12: insert a new highlighted "Exercicio 6" title+description block
# in front of the existing list (re-using the exact same wording as the old
# "Exercicio 1" block), drop the old "Exercicio 1" block, and renumber the
# remaining exercises 2,3,4,5 -> 7,8,9,10.
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$tr12 = $s12.Shapes.Item(2).TextFrame.TextRange

# New paragraph 1: "Exercicio 6" title (inherits formatting from the old
# paragraph 1 it is inserted in front of), then drop the old title text.
$tr12.Paragraphs(1, 1).InsertBefore("Exercício 6: Calcular a Média de Três Números`r") | Out-Null
$tr12.Paragraphs(2, 1).Delete()

# New paragraph 2: "Exercicio 6" description, then drop the old description.
$tr12.Paragraphs(2, 1).InsertBefore("Crie uma função anônima que recebe três números como parâmetros e retorna a média desses números.`r") | Out-Null
$tr12.Paragraphs(3, 1).Delete()

# Highlight the two new paragraphs yellow.
$tr12b = $s12.Shapes.Item(2).TextFrame2.TextRange
$tr12b.Paragraphs(1).Font.Highlight.RGB = $yellow
$tr12b.Paragraphs(2).Font.Highlight.RGB = $yellow

# Renumber "Exercicio 2" -> "Exercicio 7" (paragraph 3: title only).
$tr12.Paragraphs(3, 1).InsertBefore("Exercício 7: Calcular a Área de um Retângulo`r") | Out-Null
$tr12.Paragraphs(4, 1).Delete()

# Renumber "Exercicio 3" -> "Exercicio 8" (paragraph 5: title only).
$tr12.Paragraphs(5, 1).InsertBefore("Exercício 8: Verificar se um Número está em um Intervalo`r") | Out-Null
$tr12.Paragraphs(6, 1).Delete()

# Renumber "Exercicio 4" -> "Exercicio 9" (paragraph 7: title only).
$tr12.Paragraphs(7, 1).InsertBefore("Exercício 9: Calcular o Maior de Três Números`r") | Out-Null
$tr12.Paragraphs(8, 1).Delete()

# Renumber "Exercicio 5" -> "Exercicio 10" (paragraph 9: title only).
$tr12.Paragraphs(9, 1).InsertBefore("Exercício 10: Calcular a Soma dos Números ao Quadrado`r") | Out-Null
$tr12.Paragraphs(10, 1).Delete()

# ---------------------------------------------------------------------
# Slide 13: renumber exercises 6,7,8 -> 11,12,13 (titles only; no
# highlighting, no wording changes beyond the number).
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item(2).TextFrame.TextRange

$tr13.Paragraphs(1, 1).InsertBefore("Exercício 11: Verificar se Todos os Parâmetros são Verdadeiros `r") | Out-Null
$tr13.Paragraphs(2, 1).Delete()

$tr13.Paragraphs(3, 1).InsertBefore("Exercício 12: Calcular o Produto de Quatro Números`r") | Out-Null
$tr13.Paragraphs(4, 1).Delete()

$tr13.Paragraphs(5, 1).InsertBefore("Exercício 13: Verificar se a Soma de Dois Números é Par ou Ímpar `r") | Out-Null
$tr13.Paragraphs(6, 1).Delete()
